# "which protocols have videos" — add a "Video" column (G) marking protocols
# that have video data, plus a one-off note in H21, and expand the
# Reward_proc protocol name.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Video" header in column G
$ws.Range("G1").Value = "Video"

# Mark protocols that have video ("d") — one extra note in H21
$ws.Range("G4").Value  = "d"
$ws.Range("G9").Value  = "d"
$ws.Range("G10").Value = "d"
$ws.Range("G13").Value = "d"
$ws.Range("G17").Value = "d"
$ws.Range("G21").Value = "d"
$ws.Range("H21").Value = "preferences?"
$ws.Range("G22").Value = "i "
$ws.Range("G23").Value = "d"

# Rename protocol to reflect added NARPS data
$ws.Range("A45").Value = "Reward_proc + NARPS"

# Cosmetic: restore top-left scroll position to A1 and move the active
# selection to G41
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G41").Select()

# Header/footer font variant tweak (Normal -> Regular)
$ws.PageSetup.CenterHeader = '&"Times New Roman,Regular"&12&A'
$ws.PageSetup.CenterFooter = '&"Times New Roman,Regular"&12Page &P'
